$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19 (topic index 17), copying the formatting of row 18's A cell
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = 17

$ws.Range("B2").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.263|x15:0.737|x16:0.000|x17:0.000"
$ws.Range("B3").Value = "x1:0.000|x2:0.365|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.301|x14:0.000|x15:0.334|x16:0.000|x17:0.000"
$ws.Range("B4").Value = "x1:0.000|x2:0.000|x3:0.309|x4:0.000|x5:0.000|x6:0.148|x7:0.141|x8:0.151|x9:0.252|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B5").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.189|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.644|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.168|x17:0.000"
$ws.Range("B6").Value = "x1:0.000|x2:0.000|x3:0.284|x4:0.000|x5:0.000|x6:0.000|x7:0.209|x8:0.000|x9:0.000|x10:0.000|x11:0.154|x12:0.000|x13:0.164|x14:0.000|x15:0.000|x16:0.189|x17:0.000"
$ws.Range("B7").Value = "x1:0.139|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.861|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B8").Value = "x1:0.811|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.189|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B9").Value = "x1:0.205|x2:0.411|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.384|x17:0.000"
$ws.Range("B10").Value = "x1:0.000|x2:0.387|x3:0.000|x4:0.000|x5:0.000|x6:0.613|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B11").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:1.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B12").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:1.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B13").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:1.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B14").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:1.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B15").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.773|x9:0.000|x10:0.000|x11:0.227|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B16").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.554|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.446|x17:0.000"
$ws.Range("B17").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.610|x10:0.000|x11:0.000|x12:0.390|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B18").Value = "x1:0.000|x2:0.000|x3:0.786|x4:0.000|x5:0.000|x6:0.000|x7:0.214|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B19").Value = "x1:1.155|x2:1.163|x3:1.379|x4:0.743|x5:1.000|x6:0.761|x7:0.564|x8:1.923|x9:0.861|x10:0.833|x11:1.242|x12:0.390|x13:1.465|x14:1.263|x15:1.071|x16:1.186|x17:0.000"

Write-Host "Updated range:" $ws.UsedRange.Address()
